$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new feedback row (row 5) to the log
$row = 5

$ws.Cells.Item($row, 1).Value = "rxxx"
$ws.Cells.Item($row, 2).Value = "steve"
$ws.Cells.Item($row, 3).Value = "it worked a little"
$ws.Cells.Item($row, 4).Value = "2025-09-27 00:57:47"
